# Scheduled runner update: refresh cached Universalis price snapshots and
# recompute the dependent profit columns (H:N) across the per-job Leve
# profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
#
# Columns:
#   H currentAveragePrice        I currentAveragePriceNQ
#   J currentAveragePriceHQ      K LevePriceNQ
#   L LevePriceHQ                M LeveProfitNQ
#   N LeveProfitHQ
#
# A handful of rows gain/lose an N (or M) cell entirely where the HQ/NQ
# price feed started/stopped reporting a value for that item.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 23810356
$ws.Range("I19").Value = 28572228
$ws.Range("K19").Value = 28572228
$ws.Range("M19").Value = -28572053

$ws.Range("H107").Value = 1013.93335
$ws.Range("I107").Value = 1072.2858
$ws.Range("J107").Value = 197
$ws.Range("K107").Value = 1072.2858
$ws.Range("L107").Value = 197
$ws.Range("M107").Value = 847.7141999999999
$ws.Range("N107").Value = -4037

$ws.Range("H132").Value = 6063880
$ws.Range("I132").Value = 6454753.5
$ws.Range("J132").Value = 5340
$ws.Range("K132").Value = 19364260.5
$ws.Range("L132").Value = 16020
$ws.Range("M132").Value = -19361730.5
$ws.Range("N132").Value = -21080

$ws.Range("H137").Value = 2949.282
$ws.Range("I137").Value = 3646.7827
$ws.Range("J137").Value = 1946.625
$ws.Range("K137").Value = 10940.3481
$ws.Range("L137").Value = 5839.875
$ws.Range("M137").Value = -8390.348100000001
$ws.Range("N137").Value = -10939.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2460.5366
$ws.Range("I61").Value = 891.28
$ws.Range("K61").Value = 891.28
$ws.Range("M61").Value = -679.28

$ws.Range("H74").Value = 846
$ws.Range("I74").Value = 767.8461
$ws.Range("J74").Value = 1100
$ws.Range("K74").Value = 767.8461
$ws.Range("L74").Value = 1100
$ws.Range("M74").Value = 106.1539
$ws.Range("N74").Value = -2848

$ws.Range("H77").Value = 846
$ws.Range("I77").Value = 767.8461
$ws.Range("J77").Value = 1100
$ws.Range("K77").Value = 3839.2305
$ws.Range("L77").Value = 5500
$ws.Range("M77").Value = 528.7695000000003
$ws.Range("N77").Value = -14236

$ws.Range("H110").Value = 1100.55
$ws.Range("I110").Value = 632.1579
$ws.Range("J110").Value = 10000
$ws.Range("K110").Value = 632.1579
$ws.Range("L110").Value = 10000
$ws.Range("M110").Value = 1412.8421
$ws.Range("N110").Value = -14090

$ws.Range("H136").Value = 2460.5366
$ws.Range("I136").Value = 891.28
$ws.Range("K136").Value = 2673.84
$ws.Range("M136").Value = -123.8400000000001

$ws.Range("H139").Value = 25016.285
$ws.Range("J139").Value = 25016.285
$ws.Range("L139").Value = 25016.285
$ws.Range("N139").Value = -35296.285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4369.6665
$ws.Range("I107").Value = 2685.6667
$ws.Range("K107").Value = 2685.6667
$ws.Range("M107").Value = -765.6667000000002

$ws.Range("H134").Value = 3968.7693
$ws.Range("I134").Value = 2901.5557
$ws.Range("J134").Value = 6370
$ws.Range("K134").Value = 8704.667099999999
$ws.Range("L134").Value = 19110
$ws.Range("M134").Value = -6169.667099999999
$ws.Range("N134").Value = -24180

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 70062
$ws.Range("J42").Value = 70062
$ws.Range("L42").Value = 70062
$ws.Range("N42").Value = -71248

$ws.Range("H127").Value = 33000
$ws.Range("J127").Value = 33000
$ws.Range("L127").Value = 33000
$ws.Range("N127").Value = -42920

$ws.Range("H134").Value = 1968.7693
$ws.Range("I134").Value = 808.8
$ws.Range("J134").Value = 5835.3335
$ws.Range("K134").Value = 2426.4
$ws.Range("L134").Value = 17506.0005
$ws.Range("M134").Value = 108.6000000000004
$ws.Range("N134").Value = -22576.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1263.1666
$ws.Range("I97").Value = 589.5
$ws.Range("J97").Value = 1600
$ws.Range("K97").Value = 1768.5
$ws.Range("L97").Value = 4800
$ws.Range("M97").Value = -1272.5
$ws.Range("N97").Value = -5792

$ws.Range("H98").Value = 200
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H131").Value = 1459.7084
$ws.Range("I131").Value = 1718.8889
$ws.Range("J131").Value = 1304.2
$ws.Range("K131").Value = 5156.6667
$ws.Range("L131").Value = 3912.6
$ws.Range("M131").Value = -116.6666999999998
$ws.Range("N131").Value = -13992.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 70019
$ws.Range("J33").Value = 70019
$ws.Range("L33").Value = 70019
$ws.Range("N33").Value = -70523

$ws.Range("H38").Value = 70024
$ws.Range("J38").Value = 70024
$ws.Range("L38").Value = 70024
$ws.Range("N38").Value = -70950

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H132").Value = 4262.25
$ws.Range("I132").Value = 4926.222
$ws.Range("J132").Value = 3863.8667
$ws.Range("K132").Value = 14778.666
$ws.Range("L132").Value = 11591.6001
$ws.Range("M132").Value = -12248.666
$ws.Range("N132").Value = -16651.6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3460.3076
$ws.Range("I132").Value = 1989.6154
$ws.Range("J132").Value = 4931
$ws.Range("K132").Value = 5968.8462
$ws.Range("L132").Value = 14793
$ws.Range("M132").Value = -3438.8462
$ws.Range("N132").Value = -19853

$ws.Range("H136").Value = 1660.9231
$ws.Range("I136").Value = 1154.3448
$ws.Range("J136").Value = 3130
$ws.Range("K136").Value = 3463.0344
$ws.Range("L136").Value = 9390
$ws.Range("M136").Value = -913.0344000000005
$ws.Range("N136").Value = -14490

$ws.Range("H137").Value = 28883.572
$ws.Range("J137").Value = 28883.572
$ws.Range("L137").Value = 28883.572
$ws.Range("N137").Value = -39083.572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 70029
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 70029
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 70029
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -70663

$ws.Range("H96").Value = 810.5
$ws.Range("I96").Value = 597.3333
$ws.Range("K96").Value = 597.3333
$ws.Range("M96").Value = 775.6667

$ws.Range("H126").Value = 3707154.8
$ws.Range("I126").Value = 2745.9333
$ws.Range("J126").Value = 8337665.5
$ws.Range("K126").Value = 8237.7999
$ws.Range("L126").Value = 25012996.5
$ws.Range("M126").Value = -5767.7999
$ws.Range("N126").Value = -25017936.5

$ws.Range("H132").Value = 9757.288
$ws.Range("J132").Value = 44334.832
$ws.Range("L132").Value = 133004.496
$ws.Range("N132").Value = -138064.496

$ws.Range("H136").Value = 837.0417
$ws.Range("I136").Value = 519.5854
$ws.Range("J136").Value = 2696.4285
$ws.Range("K136").Value = 1558.7562
$ws.Range("L136").Value = 8089.2855
$ws.Range("M136").Value = 991.2437999999997
$ws.Range("N136").Value = -13189.2855

